# Auto-generated files on 2026-01-09
# Updates the HotStock_Top20 rankings table (columns A:C, rows 2-21)
# to reflect the latest "hot stock" ranking order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2"  = "金风科技"
    "A3"  = "航天电子"
    "B3"  = "岩山科技"
    "C3"  = "利欧股份"
    "A4"  = "利欧股份"
    "B4"  = "航天电子"
    "C4"  = "岩山科技"
    "A5"  = "蓝色光标"
    "B5"  = "利欧股份"
    "C5"  = "通宇通讯"
    "A6"  = "航天发展"
    "B6"  = "蓝色光标"
    "C6"  = "航天电子"
    "A7"  = "岩山科技"
    "C7"  = "航天发展"
    "A8"  = "雪人集团"
    "B8"  = "雷科防务"
    "C8"  = "雪人集团"
    "A9"  = "中国卫通"
    "B9"  = "航天发展"
    "C9"  = "雷科防务"
    "A10" = "通宇通讯"
    "B10" = "中国卫通"
    "C10" = "蓝色光标"
    "A11" = "雷科防务"
    "B11" = "南京熊猫"
    "C11" = "中国卫通"
    "A12" = "乾照光电"
    "B12" = "乾照光电"
    "C12" = "海格通信"
    "A13" = "易点天下"
    "B13" = "昆仑万维"
    "C13" = "乾照光电"
    "A14" = "昆仑万维"
    "B14" = "创新医疗"
    "C14" = "银河电子"
    "A15" = "海格通信"
    "B15" = "三花智控"
    "C15" = "华胜天成"
    "A16" = "创新医疗"
    "B16" = "通宇通讯"
    "C16" = "昆仑万维"
    "A17" = "志特新材"
    "B17" = "海格通信"
    "C17" = "创新医疗"
    "A18" = "华胜天成"
    "B18" = "中国卫星"
    "C18" = "中国卫星"
    "A19" = "三花智控"
    "B19" = "易点天下"
    "C19" = "万向钱潮"
    "B20" = "华胜天成"
    "C20" = "平潭发展"
    "A21" = "中国卫星"
    "B21" = "永鼎股份"
    "C21" = "神剑股份"
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
